{"js": "// Find the paragraph describing the C2' rotations and rewrite its text to\n// the new, condensed wording while keeping the paragraph itself in place.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldLeadIn = \"Now we have the C2\\u2019 rotations.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(oldLeadIn) === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the C2' rotations paragraph.\");\n}\n\nconst newText =\n  \"Now we have the C2\\u2019 rotations. \" +\n  \" Each rotation leaves two atoms in place, but inverts their orientation. \" +\n  \"Giving each rotation a trace of minus 2, making up for a total of -6 when all rotations are considered\" +\n  \".\";\n\ntarget.getRange().insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph describing the C2' rotations.\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"Now we have the C2\u2019 rotations.\")) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find the C2' rotations paragraph.\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n$rangeStart = $target.Range.Start\n$rangeEnd = $target.Range.End\n\n$newText = \"Now we have the C2\u2019 rotations.  Each rotation leaves two atoms in place, but inverts their orientation. Giving each rotation a trace of minus 2, making up for a total of -6 when all rotations are considered.\"\n\n# Replace across an explicit document range spanning the whole paragraph\n# (minus its trailing paragraph mark) so every run gets rewritten at once.\n$r = $d.Range($rangeStart, $rangeEnd)\n$r.Text = $newText\n"}
